# by jiankong on 1113
#
# "质控组" (QC group) is renamed to "北京组" (Beijing group) on both sheets,
# and on Sheet1 the two placeholder rows (张悦 / 卢楠, all-zero metrics) are
# removed -- the remaining 冷雪 / 屈昂 rows shift up to become rows 2-3.
# Finally Sheet1 becomes the active sheet/tab with A3 selected, and Sheet2's
# selection is reset to A2.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Drop the 张悦 row (r=2) and the 卢楠 row (now r=2 after the first delete);
# the former r=4/r=5 (冷雪/屈昂) rows shift up into r=2/r=3.
$ws1.Rows.Item(2).Delete()
$ws1.Rows.Item(2).Delete()

# Rename the group label from 质控组 to 北京组 on both sheets.
$ws1.Range("A2:A3").Value = "北京组"
$ws2.Range("A2").Value = "北京组"

# Restore the view/selection state: Sheet1 active with A3 selected,
# Sheet2's own selection reset to A2.
[void]$ws2.Range("A2").Select()
[void]$ws1.Activate()
[void]$ws1.Range("A3").Select()
